# GBFontsIssueRepro/Template_All.docx
#
# The "_GoBack" bookmark (originally wrapping "Custom fonts:" + the line
# break that follows it) is relocated so that it instead sits in the
# middle of the word "Calibri" (splitting it into "Cali" | "bri"), with
# nothing between the bookmarkStart/bookmarkEnd pair at its new home.

$d = $word.ActiveDocument

# Locate "Calibri" in the document body so we know exactly where the
# bookmark's new (empty) range belongs: right after "Cali" (4 chars in).
$findRange = $d.Content
$found = $findRange.Find.Execute("Calibri", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)

$splitPoint = $findRange.Start + 4

# Move the existing "_GoBack" bookmark to that collapsed (zero-length)
# range; Word will split the "Calibri" run into "Cali" / "bri" runs
# automatically to host the bookmark boundary.
$d.Bookmarks("_GoBack").Delete()
$target = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $target)
